$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 352; this pushes the existing rows 352-388
# down to 353-389 (and keeps their data/styles intact).
$ws.Rows("352:352").Insert()

# Populate the newly inserted row 352 with a new data record (same shape
# as the record that used to be in row 352, but with an updated date).
$ws.Range("A352").Value = 4
$ws.Range("B352").Value = 'Feria Lagunitas de Puerto Montt'
$ws.Range("C352").Value = 'Los Lagos'
$ws.Range("D352").Value = 44946
$ws.Range("E352").Value = 10
$ws.Range("F352").Value = 100112037
$ws.Range("G352").Value = 'Cebollín'
$ws.Range("H352").Value = 'Sin especificar'
$ws.Range("I352").Value = 'Primera'
$ws.Range("J352").Value = 180
$ws.Range("K352").Value = 6000
$ws.Range("L352").Value = 6000
$ws.Range("M352").Value = 6000
$ws.Range("N352").Value = '$/paquete 36 unidades'
$ws.Range("O352").Value = 'Región Metropolitana'
$ws.Range("P352").Value = 167
$ws.Range("Q352").Value = 36
$ws.Range("R352").Value = 'Hortaliza'
